# chore: update Sheets via scheduled runner
# Refreshes market-price-derived columns (H..N) for a set of leve rows
# across all 8 crafter sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 466.66666
$ws.Cells.Item(33, 9).Value = 407.2857
$ws.Cells.Item(33, 10).Value = 674.5
$ws.Cells.Item(33, 11).Value = 407.2857
$ws.Cells.Item(33, 12).Value = 674.5
$ws.Cells.Item(33, 13).Value = -178.2857
$ws.Cells.Item(33, 14).Value = -1132.5

$ws.Cells.Item(55, 8).Value = 923.5
$ws.Cells.Item(55, 9).Value = 602.5
$ws.Cells.Item(55, 11).Value = 602.5
$ws.Cells.Item(55, 13).Value = -388.5

$ws.Cells.Item(116, 8).Value = 32993.35
$ws.Cells.Item(116, 9).Value = 40237
$ws.Cells.Item(116, 10).Value = 26554.555
$ws.Cells.Item(116, 11).Value = 40237
$ws.Cells.Item(116, 12).Value = 26554.555
$ws.Cells.Item(116, 13).Value = -36795
$ws.Cells.Item(116, 14).Value = -33438.555

$ws.Cells.Item(132, 8).Value = 1759.1464
$ws.Cells.Item(132, 9).Value = 1624.9459
$ws.Cells.Item(132, 11).Value = 4874.8377
$ws.Cells.Item(132, 13).Value = -2344.8377

$ws.Cells.Item(141, 8).Value = 4217.7188
$ws.Cells.Item(141, 9).Value = 3180.9656
$ws.Cells.Item(141, 10).Value = 14239.667
$ws.Cells.Item(141, 11).Value = 9542.8968
$ws.Cells.Item(141, 12).Value = 42719.001
$ws.Cells.Item(141, 13).Value = -4362.8968
$ws.Cells.Item(141, 14).Value = -53079.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(31, 8).Value = 29078.166
$ws.Cells.Item(31, 9).Value = 6117.5
$ws.Cells.Item(31, 11).Value = 6117.5
$ws.Cells.Item(31, 13).Value = -5823.5

$ws.Cells.Item(61, 8).Value = 6003741.5
$ws.Cells.Item(61, 9).Value = 7146202
$ws.Cells.Item(61, 11).Value = 7146202
$ws.Cells.Item(61, 13).Value = -7145990

$ws.Cells.Item(74, 8).Value = 2343.4614
$ws.Cells.Item(74, 9).Value = 1507.25
$ws.Cells.Item(74, 11).Value = 1507.25
$ws.Cells.Item(74, 13).Value = -633.25

$ws.Cells.Item(77, 8).Value = 2343.4614
$ws.Cells.Item(77, 9).Value = 1507.25
$ws.Cells.Item(77, 11).Value = 7536.25
$ws.Cells.Item(77, 13).Value = -3168.25

$ws.Cells.Item(132, 8).Value = 5558870
$ws.Cells.Item(132, 9).Value = 3564.7693
$ws.Cells.Item(132, 11).Value = 10694.3079
$ws.Cells.Item(132, 13).Value = -8164.3079

$ws.Cells.Item(136, 8).Value = 6003741.5
$ws.Cells.Item(136, 9).Value = 7146202
$ws.Cells.Item(136, 11).Value = 21438606
$ws.Cells.Item(136, 13).Value = -21436056

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3069.8
$ws.Cells.Item(86, 9).Value = 2155.2942
$ws.Cells.Item(86, 10).Value = 8252
$ws.Cells.Item(86, 11).Value = 2155.2942
$ws.Cells.Item(86, 12).Value = 8252
$ws.Cells.Item(86, 13).Value = -1032.2942
$ws.Cells.Item(86, 14).Value = -10498

$ws.Cells.Item(89, 8).Value = 3069.8
$ws.Cells.Item(89, 9).Value = 2155.2942
$ws.Cells.Item(89, 10).Value = 8252
$ws.Cells.Item(89, 11).Value = 10776.471
$ws.Cells.Item(89, 12).Value = 41260
$ws.Cells.Item(89, 13).Value = -5160.471
$ws.Cells.Item(89, 14).Value = -52492

$ws.Cells.Item(94, 8).Value = 1898.75
$ws.Cells.Item(94, 9).Value = 2458.7144
$ws.Cells.Item(94, 10).Value = 1463.2222
$ws.Cells.Item(94, 11).Value = 2458.7144
$ws.Cells.Item(94, 12).Value = 1463.2222
$ws.Cells.Item(94, 13).Value = -2007.7144
$ws.Cells.Item(94, 14).Value = -2365.2222

$ws.Cells.Item(134, 8).Value = 5557845
$ws.Cells.Item(134, 9).Value = 2324.625
$ws.Cells.Item(134, 11).Value = 6973.875
$ws.Cells.Item(134, 13).Value = -4438.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 5000
$ws.Cells.Item(4, 10).Value = 5000
$ws.Cells.Item(4, 12).Value = 5000
$ws.Cells.Item(4, 14).Value = -5224

$ws.Cells.Item(31, 8).Value = 34485572
$ws.Cells.Item(31, 10).Value = 4201.857
$ws.Cells.Item(31, 12).Value = 4201.857
$ws.Cells.Item(31, 14).Value = -4791.857

$ws.Cells.Item(34, 8).Value = 34485572
$ws.Cells.Item(34, 10).Value = 4201.857
$ws.Cells.Item(34, 12).Value = 4201.857
$ws.Cells.Item(34, 14).Value = -4605.857

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(32, 8).Value = 2515332
$ws.Cells.Item(32, 10).Value = 3020442.8
$ws.Cells.Item(32, 12).Value = 9061328.399999999
$ws.Cells.Item(32, 14).Value = -9061894.399999999

$ws.Cells.Item(70, 8).Value = 21778.467
$ws.Cells.Item(70, 9).Value = 20012
$ws.Cells.Item(70, 11).Value = 60036
$ws.Cells.Item(70, 13).Value = -59721

$ws.Cells.Item(73, 8).Value = 21778.467
$ws.Cells.Item(73, 9).Value = 20012
$ws.Cells.Item(73, 11).Value = 60036
$ws.Cells.Item(73, 13).Value = -58944

$ws.Cells.Item(75, 8).Value = 11979.889
$ws.Cells.Item(75, 9).Value = 8500
$ws.Cells.Item(75, 10).Value = 12414.875
$ws.Cells.Item(75, 11).Value = 25500
$ws.Cells.Item(75, 12).Value = 37244.625
$ws.Cells.Item(75, 13).Value = -24502
$ws.Cells.Item(75, 14).Value = -39240.625

$ws.Cells.Item(78, 8).Value = 11979.889
$ws.Cells.Item(78, 9).Value = 8500
$ws.Cells.Item(78, 10).Value = 12414.875
$ws.Cells.Item(78, 11).Value = 76500
$ws.Cells.Item(78, 12).Value = 111733.875
$ws.Cells.Item(78, 13).Value = -71508
$ws.Cells.Item(78, 14).Value = -121717.875

$ws.Cells.Item(117, 8).Value = 5737.385
$ws.Cells.Item(117, 9).Value = 980
$ws.Cells.Item(117, 11).Value = 2940
$ws.Cells.Item(117, 13).Value = 502

$ws.Cells.Item(121, 8).Value = 4730.4443
$ws.Cells.Item(121, 9).Value = 499.5
$ws.Cells.Item(121, 10).Value = 5259.3125
$ws.Cells.Item(121, 11).Value = 1498.5
$ws.Cells.Item(121, 12).Value = 15777.9375
$ws.Cells.Item(121, 13).Value = -188.5
$ws.Cells.Item(121, 14).Value = -18397.9375

$ws.Cells.Item(122, 8).Value = 41664.25
$ws.Cells.Item(122, 9).Value = 66367
$ws.Cells.Item(122, 11).Value = 597303
$ws.Cells.Item(122, 13).Value = -594853

$ws.Cells.Item(129, 8).Value = 9432.117
$ws.Cells.Item(129, 10).Value = 9141.154
$ws.Cells.Item(129, 12).Value = 27423.462
$ws.Cells.Item(129, 14).Value = -37423.462

$ws.Cells.Item(139, 8).Value = 5062.5
$ws.Cells.Item(139, 9).Value = 1233.0714
$ws.Cells.Item(139, 10).Value = 10423.7
$ws.Cells.Item(139, 11).Value = 3699.2142
$ws.Cells.Item(139, 12).Value = 31271.1
$ws.Cells.Item(139, 13).Value = 1440.7858
$ws.Cells.Item(139, 14).Value = -41551.10000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 14).ClearContents()

$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 14).ClearContents()

$ws.Cells.Item(126, 8).Value = 5029.25
$ws.Cells.Item(126, 9).Value = 4264.6924
$ws.Cells.Item(126, 11).Value = 12794.0772
$ws.Cells.Item(126, 13).Value = -10324.0772

$ws.Cells.Item(128, 8).Value = 80798.39999999999
$ws.Cells.Item(128, 10).Value = 80798.39999999999
$ws.Cells.Item(128, 12).Value = 80798.39999999999
$ws.Cells.Item(128, 14).Value = -90758.39999999999

$ws.Cells.Item(129, 8).Value = 99999
$ws.Cells.Item(129, 10).Value = 99999
$ws.Cells.Item(129, 12).Value = 99999
$ws.Cells.Item(129, 14).Value = -109999

$ws.Cells.Item(130, 8).Value = 99999
$ws.Cells.Item(130, 10).Value = 99999
$ws.Cells.Item(130, 12).Value = 99999
$ws.Cells.Item(130, 14).Value = -110039

$ws.Cells.Item(132, 8).Value = 9097066
$ws.Cells.Item(132, 9).Value = 6966.75
$ws.Cells.Item(132, 11).Value = 20900.25
$ws.Cells.Item(132, 13).Value = -18370.25

$ws.Cells.Item(134, 8).Value = 51122.445
$ws.Cells.Item(134, 10).Value = 51122.445
$ws.Cells.Item(134, 12).Value = 153367.335
$ws.Cells.Item(134, 14).Value = -158437.335

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 9269.261
$ws.Cells.Item(7, 9).Value = 9210.833000000001
$ws.Cells.Item(7, 10).Value = 9479.6
$ws.Cells.Item(7, 11).Value = 9210.833000000001
$ws.Cells.Item(7, 12).Value = 9479.6
$ws.Cells.Item(7, 13).Value = -9098.833000000001
$ws.Cells.Item(7, 14).Value = -9703.6

$ws.Cells.Item(22, 8).Value = 2828.6667
$ws.Cells.Item(22, 9).Value = 1730
$ws.Cells.Item(22, 11).Value = 1730
$ws.Cells.Item(22, 13).Value = -1435

$ws.Cells.Item(27, 8).Value = 2828.6667
$ws.Cells.Item(27, 9).Value = 1730
$ws.Cells.Item(27, 11).Value = 1730
$ws.Cells.Item(27, 13).Value = -1623

$ws.Cells.Item(126, 8).Value = 9269.261
$ws.Cells.Item(126, 9).Value = 9210.833000000001
$ws.Cells.Item(126, 10).Value = 9479.6
$ws.Cells.Item(126, 11).Value = 27632.499
$ws.Cells.Item(126, 12).Value = 28438.8
$ws.Cells.Item(126, 13).Value = -25162.499
$ws.Cells.Item(126, 14).Value = -33378.8

$ws.Cells.Item(132, 8).Value = 4677.619
$ws.Cells.Item(132, 9).Value = 2856.3
$ws.Cells.Item(132, 10).Value = 6333.364
$ws.Cells.Item(132, 11).Value = 8568.900000000001
$ws.Cells.Item(132, 12).Value = 19000.092
$ws.Cells.Item(132, 13).Value = -6038.900000000001
$ws.Cells.Item(132, 14).Value = -24060.092

$ws.Cells.Item(136, 8).Value = 5164.2964
$ws.Cells.Item(136, 9).Value = 4271.222
$ws.Cells.Item(136, 10).Value = 6950.4443
$ws.Cells.Item(136, 11).Value = 12813.666
$ws.Cells.Item(136, 12).Value = 20851.3329
$ws.Cells.Item(136, 13).Value = -10263.666
$ws.Cells.Item(136, 14).Value = -25951.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1570
$ws.Cells.Item(81, 9).Value = 1582.8572
$ws.Cells.Item(81, 11).Value = 3165.7144
$ws.Cells.Item(81, 13).Value = -2104.7144

$ws.Cells.Item(84, 8).Value = 1570
$ws.Cells.Item(84, 9).Value = 1582.8572
$ws.Cells.Item(84, 11).Value = 15828.572
$ws.Cells.Item(84, 13).Value = -10524.572

$ws.Cells.Item(132, 8).Value = 287653.4
$ws.Cells.Item(132, 9).Value = 1867.48
$ws.Cells.Item(132, 10).Value = 1002118.2
$ws.Cells.Item(132, 11).Value = 5602.440000000001
$ws.Cells.Item(132, 12).Value = 3006354.6
$ws.Cells.Item(132, 13).Value = -3072.440000000001
$ws.Cells.Item(132, 14).Value = -3011414.6

$ws.Cells.Item(136, 8).Value = 287110
$ws.Cells.Item(136, 9).Value = 1180.5416
$ws.Cells.Item(136, 10).Value = 910956.0600000001
$ws.Cells.Item(136, 11).Value = 3541.6248
$ws.Cells.Item(136, 12).Value = 2732868.18
$ws.Cells.Item(136, 13).Value = -991.6248000000001
$ws.Cells.Item(136, 14).Value = -2737968.18
